# Scheduled data-refresh: update market-price & profit columns (H-N)
# for the affected Leve rows across each job sheet in the workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5000
$ws.Range("J18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("N18").Value = -5568

$ws.Range("H86").Value = 1499.5
$ws.Range("I86").Value = 999.5
$ws.Range("J86").Value = 1999.5
$ws.Range("K86").Value = 999.5
$ws.Range("L86").Value = 1999.5
$ws.Range("M86").Value = 123.5
$ws.Range("N86").Value = -4245.5

$ws.Range("H89").Value = 1499.5
$ws.Range("I89").Value = 999.5
$ws.Range("J89").Value = 1999.5
$ws.Range("K89").Value = 4997.5
$ws.Range("L89").Value = 9997.5
$ws.Range("M89").Value = 618.5
$ws.Range("N89").Value = -21229.5

$ws.Range("H112").Value = 6246.25
$ws.Range("J112").Value = 6246.25
$ws.Range("L112").Value = 18738.75
$ws.Range("N112").Value = -20954.75

$ws.Range("H127").Value = 745.44446
$ws.Range("I127").Value = 686.5
$ws.Range("K127").Value = 2059.5
$ws.Range("M127").Value = 2900.5

$ws.Range("H137").Value = 3190.5925
$ws.Range("I137").Value = 2274.6
$ws.Range("J137").Value = 3729.4119
$ws.Range("K137").Value = 6823.799999999999
$ws.Range("L137").Value = 11188.2357
$ws.Range("M137").Value = -4273.799999999999
$ws.Range("N137").Value = -16288.2357

$ws.Range("H138").Value = 4417.174
$ws.Range("I138").Value = 4090.75
$ws.Range("J138").Value = 4591.2666
$ws.Range("K138").Value = 12272.25
$ws.Range("L138").Value = 13773.7998
$ws.Range("M138").Value = -7132.25
$ws.Range("N138").Value = -24053.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3994
$ws.Range("I61").Value = 3994
$ws.Range("K61").Value = 3994
$ws.Range("M61").Value = -3782

$ws.Range("H74").Value = 1724.875
$ws.Range("I74").Value = 1724.875
$ws.Range("K74").Value = 1724.875
$ws.Range("M74").Value = -850.875

$ws.Range("H77").Value = 1724.875
$ws.Range("I77").Value = 1724.875
$ws.Range("K77").Value = 8624.375
$ws.Range("M77").Value = -4256.375

$ws.Range("H122").Value = 1004.3333
$ws.Range("I122").Value = 999.5
$ws.Range("K122").Value = 2998.5
$ws.Range("M122").Value = -548.5

$ws.Range("H136").Value = 3994
$ws.Range("I136").Value = 3994
$ws.Range("K136").Value = 11982
$ws.Range("M136").Value = -9432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 874.6
$ws.Range("I80").Value = 789.3333
$ws.Range("J80").Value = 1002.5
$ws.Range("K80").Value = 789.3333
$ws.Range("L80").Value = 1002.5
$ws.Range("M80").Value = 208.6667
$ws.Range("N80").Value = -2998.5

$ws.Range("H83").Value = 874.6
$ws.Range("I83").Value = 789.3333
$ws.Range("J83").Value = 1002.5
$ws.Range("K83").Value = 3946.6665
$ws.Range("L83").Value = 5012.5
$ws.Range("M83").Value = 1045.3335
$ws.Range("N83").Value = -14996.5

$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 5000
$ws.Range("K86").Value = 5000
$ws.Range("M86").Value = -3877

$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 5000
$ws.Range("K89").Value = 25000
$ws.Range("M89").Value = -19384

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35569.715
$ws.Range("I31").Value = 16498.666
$ws.Range("K31").Value = 16498.666
$ws.Range("M31").Value = -16203.666

$ws.Range("H34").Value = 35569.715
$ws.Range("I34").Value = 16498.666
$ws.Range("K34").Value = 16498.666
$ws.Range("M34").Value = -16296.666

$ws.Range("H132").Value = 1937.8182
$ws.Range("I132").Value = 1757.3334
$ws.Range("K132").Value = 5272.0002
$ws.Range("M132").Value = -2742.0002

$ws.Range("H134").Value = 1126.6
$ws.Range("I134").Value = 1085.1111
$ws.Range("K134").Value = 3255.3333
$ws.Range("M134").Value = -720.3333000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 111392.11
$ws.Range("I4").Value = 218.57143
$ws.Range("J4").Value = 500499.5
$ws.Range("K4").Value = 655.71429
$ws.Range("L4").Value = 1501498.5
$ws.Range("M4").Value = -543.71429
$ws.Range("N4").Value = -1501722.5

$ws.Range("H13").Value = 12
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 12
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 36
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -372

$ws.Range("H104").Value = 1995
$ws.Range("I104").Value = 1995
$ws.Range("K104").Value = 5985
$ws.Range("M104").Value = -3364

$ws.Range("H110").Value = 3000
$ws.Range("I110").Value = 3000
$ws.Range("K110").Value = 9000
$ws.Range("M110").Value = -4910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5553.4287
$ws.Range("I102").Value = 5519.077
$ws.Range("K102").Value = 5519.077
$ws.Range("M102").Value = -3897.077

$ws.Range("H126").Value = 4705.8667
$ws.Range("I126").Value = 5196.75
$ws.Range("K126").Value = 15590.25
$ws.Range("M126").Value = -13120.25

$ws.Range("H132").Value = 4518.625
$ws.Range("I132").Value = 4358.1665
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 13074.4995
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -10544.4995
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 534.4
$ws.Range("I16").Value = 534.4
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 534.4
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -364.4
$ws.Range("N16").ClearContents()

$ws.Range("H22").Value = 774.75
$ws.Range("J22").Value = 849.5
$ws.Range("L22").Value = 849.5
$ws.Range("N22").Value = -1439.5

$ws.Range("H27").Value = 774.75
$ws.Range("J27").Value = 849.5
$ws.Range("L27").Value = 849.5
$ws.Range("N27").Value = -1063.5

$ws.Range("H40").Value = 4251
$ws.Range("I40").Value = 4251
$ws.Range("K40").Value = 4251
$ws.Range("M40").Value = -4115

$ws.Range("H82").Value = 1999.625
$ws.Range("I82").Value = 1999.625
$ws.Range("K82").Value = 1999.625
$ws.Range("M82").Value = -1638.625

$ws.Range("H85").Value = 1999.625
$ws.Range("I85").Value = 1999.625
$ws.Range("K85").Value = 1999.625
$ws.Range("M85").Value = -751.625

$ws.Range("H93").Value = 1997.375
$ws.Range("I93").Value = 1997.375
$ws.Range("K93").Value = 1997.375
$ws.Range("M93").Value = -749.375

$ws.Range("H132").Value = 22647.1
$ws.Range("I132").Value = 21933.875
$ws.Range("J132").Value = 25500
$ws.Range("K132").Value = 65801.625
$ws.Range("L132").Value = 76500
$ws.Range("M132").Value = -63271.625
$ws.Range("N132").Value = -81560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2030.375
$ws.Range("I126").Value = 1040.6666
$ws.Range("K126").Value = 3121.9998
$ws.Range("M126").Value = -651.9998000000001

$ws.Range("H132").Value = 10589.3
$ws.Range("I132").Value = 8379.4
$ws.Range("J132").Value = 12799.2
$ws.Range("K132").Value = 25138.2
$ws.Range("L132").Value = 38397.60000000001
$ws.Range("M132").Value = -22608.2
$ws.Range("N132").Value = -43457.60000000001

$ws.Range("H137").Value = 100699
$ws.Range("J137").Value = 100699
$ws.Range("L137").Value = 100699
$ws.Range("N137").Value = -110899
